$d = $word.ActiveDocument

# 1. Remove the hidden "_GoBack" bookmark from the title paragraph; it will
#    be re-created at the end of the new "Q. " paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Split "Neural Networks: Learning" into "Neural Networks: " + "Learning"
#    (the word "Learning" becomes its own run).
$d.Content.Find.Execute(": Learning", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ": ", 2) | Out-Null

$titlePara = $d.Paragraphs(1)
$insPos = $titlePara.Range.End - 1
$d.Range($insPos, $insPos).InsertAfter("Learning") | Out-Null

# 3. Drop the quiz table entirely.
$d.Tables(1).Delete()

# 4. Collapse the two trailing empty paragraphs (that used to follow the
#    table) down to a single paragraph, keeping the empty paragraph that
#    used to sit right before the table.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Range($lastPara.Range.Start - 1, $lastPara.Range.End).Delete() | Out-Null

# 5. Turn the now-last paragraph into "Q. " with English (US) paragraph-mark
#    formatting and re-attach the "_GoBack" bookmark at its end.
$qPara = $d.Paragraphs($d.Paragraphs.Count)
$qPara.Range.LanguageID = "en-US"
$qPara.Range.InsertAfter("Q. ") | Out-Null

$qEnd = $qPara.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($qEnd, $qEnd)) | Out-Null
